$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

# Status text: "Handed back: in sync with en-US" -> "Ready for handoff"
$overview.Range("E2").Value = "Ready for handoff"
$overview.Range("F2").Value = "Ready for handoff"
$zhcn.Range("C2").Value = "Ready for handoff"
$dede.Range("C2").Value = "Ready for handoff"

# Datetimes refreshed for the new handoff report
$overview.Range("G2").Value = "2016-10-27 10:06:35"
$dede.Range("H2").Value = "2016-10-27 10:06:35"
$zhcn.Range("H2").Value = "2016-10-27 10:06:22"

# The shorter status text lets these columns shrink to fit (AutoFit-style
# resize down from the old ~30-char-wide column to fit "Ready for handoff")
$overview.Columns.Item(5).ColumnWidth = 16.33
$overview.Columns.Item(6).ColumnWidth = 16.33
$zhcn.Columns.Item(3).ColumnWidth = 16.33
$dede.Columns.Item(3).ColumnWidth = 16.33

